# fix: ensure elective courses are scheduled in same time slots for both sections A and B
#
# 1. Update Section_A / Section_B timetables so the elective slots line up:
#      E5: Free              -> DS401 (Elective)
#      F5: EC460 (Elective)  -> EC456 (Elective)
#      E6: DS460 (Elective)  -> DS456 (Elective)
# 2. Update Course_Summary with the renumbered/renamed elective courses and
#    their new instructors.
# 3. Add a new "Elective_Coordination" sheet listing each elective course
#    with its day/time slot/slot-name and the sections sharing that slot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Section_A and Section_B timetable fixes
# ---------------------------------------------------------------------------
foreach ($sheetName in @("Section_A", "Section_B")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("E5").Value = "DS401 (Elective)"
    $ws.Range("F5").Value = "EC456 (Elective)"
    $ws.Range("E6").Value = "DS456 (Elective)"
}

# ---------------------------------------------------------------------------
# 2. Course_Summary updates
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Course_Summary")

# Row 2: DS460 Security Analytics -> DS456 Cybersecurity Techniques
$summary.Range("A2").Value = "DS456"
$summary.Range("B2").Value = "Cybersecurity Techniques"
$summary.Range("F2").Value = "Dr. Rajendra Hegadi"

# Row 3: EC460 Deep Learning -> EC456 Reinforcement Learning
$summary.Range("A3").Value = "EC456"
$summary.Range("B3").Value = "Reinforcement Learning"
$summary.Range("F3").Value = "Dr. Divyajyothi"

# Row 4: DS405 Bioinformatics -> DS401 Health Care Data Analytics
$summary.Range("A4").Value = "DS401"
$summary.Range("B4").Value = "Health Care Data Analytics"
$summary.Range("F4").Value = "Dr. Girish G N"

# ---------------------------------------------------------------------------
# 3. New Elective_Coordination sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$coord = $wb.Worksheets.Add($null, $lastSheet)
$coord.Name = "Elective_Coordination"

$headers = @("Elective Course", "Day", "Time Slot", "Slot Name", "Sections")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $coord.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @("DS456", "Thu", "15:30-17:00", "Elective_Slot_1", "A & B (Common Slot)"),
    @("EC456", "Fri", "14:00-15:30", "Elective_Slot_2", "A & B (Common Slot)"),
    @("DS401", "Thu", "14:00-15:30", "Elective_Slot_3", "A & B (Common Slot)")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $coord.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}
